# Auto-generated: add rows 182-189 to each of the 4 sheets, matching the target diff.
$wb = $excel.ActiveWorkbook

# --- Fix sheet1 (DE_LFT_#1) row 181, column A: corrected timestamp ---
$wsFix = $wb.Worksheets.Item("DE_LFT_#1")
$wsFix.Range("A181").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsFix.Range("A181").Value = [double]"45967.43369212963"

# --- Sheet 1 (DE_LFT_#1): append rows 182-189 ---
$ws1 = $wb.Worksheets.Item("DE_LFT_#1")
$rows1 = @()
$rows1 += , @(182, "45968.43420138889", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xD1", "0x14", 380, "7.59863127514711e+23", 200, 14)
$rows1 += , @(183, "45969.43471064815", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,1xCC", "0x14", 380, "7.59863127514711e+23", 200, 14)
$rows1 += , @(184, "45970.43521990741", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xD2", "0x14", 380, "7.59863127514711e+23", 196, 14)
$rows1 += , @(185, "45971.43572916667", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,2xCC", "0x14", 380, "7.59863127514711e+23", 192, 14)
$rows1 += , @(186, "45972.43623842593", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xD3", "0x14", 380, "7.59863127514711e+23", 188, 14)
$rows1 += , @(187, "45973.43674768518", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,3xCC", "0x14", 380, "7.59863127514711e+23", 188, 14)
$rows1 += , @(188, "45974.43725694445", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,0xD4", "0x14", 380, "7.59863127514711e+23", 184, 14)
$rows1 += , @(189, "45975.4377662037", "0x01,0x7c", "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,", "0x00,4xCC", "0x14", 380, "7.59863127514711e+23", 180, 14)
foreach ($r in $rows1) {
    $rn = $r[0]
    $ws1.Cells.Item($rn, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws1.Cells.Item($rn, 1).Value = [double]$r[1]
    $ws1.Cells.Item($rn, 2).Value = $r[2]
    $ws1.Cells.Item($rn, 3).Value = $r[3]
    $ws1.Cells.Item($rn, 4).Value = $r[4]
    $ws1.Cells.Item($rn, 5).Value = $r[5]
    $ws1.Cells.Item($rn, 6).Value = [double]$r[6]
    $ws1.Cells.Item($rn, 7).Value = [double]$r[7]
    $ws1.Cells.Item($rn, 8).Value = [double]$r[8]
    $ws1.Cells.Item($rn, 9).Value = [double]$r[9]
}

# --- Sheet 2 (DE_LFT_#2): append rows 182-189 ---
$ws2 = $wb.Worksheets.Item("DE_LFT_#2")
$rows2 = @()
$rows2 += , @(182, "45968.43420138889", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD0", "0xe", 380, "5.68432987514711e+23", 208, 14)
$rows2 += , @(183, "45969.43471064815", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD4", "0xe", 380, "5.68432987514711e+23", 204, 14)
$rows2 += , @(184, "45970.43521990741", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD8", "0xe", 380, "5.68432987514711e+23", 204, 14)
$rows2 += , @(185, "45971.43572916667", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD12", "0xe", 380, "5.68432987514711e+23", 200, 14)
$rows2 += , @(186, "45972.43623842593", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD16", "0xe", 380, "5.68432987514711e+23", 196, 14)
$rows2 += , @(187, "45973.43674768518", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD20", "0xe", 380, "5.68432987514711e+23", 196, 14)
$rows2 += , @(188, "45974.43725694445", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD24", "0xe", 380, "5.68432987514711e+23", 192, 14)
$rows2 += , @(189, "45975.4377662037", "0x01,0x7c", "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,", "0x00,0xD28", "0xe", 380, "5.68432987514711e+23", 188, 14)
foreach ($r in $rows2) {
    $rn = $r[0]
    $ws2.Cells.Item($rn, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws2.Cells.Item($rn, 1).Value = [double]$r[1]
    $ws2.Cells.Item($rn, 2).Value = $r[2]
    $ws2.Cells.Item($rn, 3).Value = $r[3]
    $ws2.Cells.Item($rn, 4).Value = $r[4]
    $ws2.Cells.Item($rn, 5).Value = $r[5]
    $ws2.Cells.Item($rn, 6).Value = [double]$r[6]
    $ws2.Cells.Item($rn, 7).Value = [double]$r[7]
    $ws2.Cells.Item($rn, 8).Value = [double]$r[8]
    $ws2.Cells.Item($rn, 9).Value = [double]$r[9]
}

# --- Sheet 3 (DE_PLT_#1): append rows 182-189 ---
$ws3 = $wb.Worksheets.Item("DE_PLT_#1")
$rows3 = @()
$rows3 += , @(182, "45968.43420138889", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x65", "0x7", 130, "5.68631262647114e+23", 101, 7)
$rows3 += , @(183, "45969.43471064815", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x64", "0x7", 130, "5.68631262647114e+23", 94, 7)
$rows3 += , @(184, "45970.43521990741", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x63", "0x7", 130, "5.68631262647114e+23", 94, 7)
$rows3 += , @(185, "45971.43572916667", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x62", "0x7", 130, "5.68631262647114e+23", 93, 7)
$rows3 += , @(186, "45972.43623842593", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x61", "0x7", 130, "5.68631262647114e+23", 93, 7)
$rows3 += , @(187, "45973.43674768518", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x60", "0x7", 130, "5.68631262647114e+23", 90, 7)
$rows3 += , @(188, "45974.43725694445", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x59", "0x7", 130, "5.68631262647114e+23", 90, 7)
$rows3 += , @(189, "45975.4377662037", "0x00,0x82", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,", "0x00,0x58", "0x7", 130, "5.68631262647114e+23", 88, 7)
foreach ($r in $rows3) {
    $rn = $r[0]
    $ws3.Cells.Item($rn, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws3.Cells.Item($rn, 1).Value = [double]$r[1]
    $ws3.Cells.Item($rn, 2).Value = $r[2]
    $ws3.Cells.Item($rn, 3).Value = $r[3]
    $ws3.Cells.Item($rn, 4).Value = $r[4]
    $ws3.Cells.Item($rn, 5).Value = $r[5]
    $ws3.Cells.Item($rn, 6).Value = [double]$r[6]
    $ws3.Cells.Item($rn, 7).Value = [double]$r[7]
    $ws3.Cells.Item($rn, 8).Value = [double]$r[8]
    $ws3.Cells.Item($rn, 9).Value = [double]$r[9]
}

# --- Sheet 4 (DE_PLT_#2): append rows 182-189 ---
$ws4 = $wb.Worksheets.Item("DE_PLT_#2")
$rows4 = @()
$rows4 += , @(182, "45968.43420138889", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x62", "0x3", 130, "9.85046333984776e+23", 99, 3)
$rows4 += , @(183, "45969.43471064815", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x61", "0x3", 130, "9.85046333984776e+23", 98, 3)
$rows4 += , @(184, "45970.43521990741", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x60", "0x3", 130, "9.85046333984776e+23", 98, 3)
$rows4 += , @(185, "45971.43572916667", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x59", "0x3", 130, "9.85046333984776e+23", 97, 3)
$rows4 += , @(186, "45972.43623842593", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x58", "0x3", 130, "9.85046333984776e+23", 97, 3)
$rows4 += , @(187, "45973.43674768518", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x57", "0x3", 130, "9.85046333984776e+23", 96, 3)
$rows4 += , @(188, "45974.43725694445", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x56", "0x3", 130, "9.85046333984776e+23", 96, 3)
$rows4 += , @(189, "45975.4377662037", "0x00,0x82", "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,", "0x00,0x55", "0x3", 130, "9.85046333984776e+23", 95, 3)
foreach ($r in $rows4) {
    $rn = $r[0]
    $ws4.Cells.Item($rn, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws4.Cells.Item($rn, 1).Value = [double]$r[1]
    $ws4.Cells.Item($rn, 2).Value = $r[2]
    $ws4.Cells.Item($rn, 3).Value = $r[3]
    $ws4.Cells.Item($rn, 4).Value = $r[4]
    $ws4.Cells.Item($rn, 5).Value = $r[5]
    $ws4.Cells.Item($rn, 6).Value = [double]$r[6]
    $ws4.Cells.Item($rn, 7).Value = [double]$r[7]
    $ws4.Cells.Item($rn, 8).Value = [double]$r[8]
    $ws4.Cells.Item($rn, 9).Value = [double]$r[9]
}

Write-Host "Done adding rows."
